# Updates cryptos list (prices + 1h volume %) per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "64.553.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "'" + "3.142.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'" + "572.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("D6").Value = "'" + "148.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'" + "3.140.38"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("D9").Value = "'" + "0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("D11").Value = "'" + "6.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").Value = "'" + "0.0000263"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.23%  "
$ws.Range("D14").Value = "'" + "36.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").Value = "'" + "3.654.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").Value = "'" + "64.755.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "'" + "3.143.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").Value = "'" + "7.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "'" + "501.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("D21").Value = "'" + "14.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("D22").Value = "'" + "0.714"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").Value = "'" + "15.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "'" + "83.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "'" + "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("D27").Value = "'" + "8.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.30%  "
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("D29").Value = "'" + "2.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("D30").Value = "'" + "2.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.66%  "
$ws.Range("D31").Value = "'" + "27.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("E34").Value = "  +3.19%  "
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").Value = "'" + "54.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.70%  "
$ws.Range("E37").Value = "  +6.84%  "
$ws.Range("D38").Value = "'" + "466.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "'" + "2.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.91%  "
$ws.Range("D41").Value = "'" + "8.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'" + "3.032.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.71%  "
$ws.Range("D43").Value = "'" + "0.116"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("B44").Value = "'" + "TheGraph"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'" + "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'" + "0.282"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("B45").Value = "'" + "Fetch.AI"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'" + "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'" + "2.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.01%  "
$ws.Range("D46").Value = "'" + "28.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").Value = "'" + "0.0₃0578"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.13%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").Value = "'" + "118.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.84%  "
